# "fixed the expense link bug"
# A new expense row ("Party", 2000, 28/5/2025) was missing from the table;
# insert it as the new row 2 (right after the header), which pushes the
# existing Cafe/movie/travel/clothing/rent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 2 ("Cafe"), shifting
# Cafe/movie/travel/clothing/rent down to rows 3-7.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the missing expense entry.
$ws.Range("A2").Value = "Party"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = "28/5/2025"
